$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the format of the existing header cell (H1) onto the new header cells
# so they share the same bold/centered/bordered style instead of creating
# a brand-new style entry.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Data values for the new columns I (I0) and J (IF), rows 2-17
$dataI = @(6, 5, 8, 6, 8, 6, 5, 8, 5, 7, 7, 5, 7, 7, 5, 7)
$dataJ = @(7, 6, 8, 7, 8, 7, 6, 9, 6, 7, 8, 5, 9, 7, 6, 8)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
